# Auto-generated edit script applying numeric value updates to the Raiden_Profits workbook
# (sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 518.7143
$ws.Range("I33").Value = 361.0625
$ws.Range("K33").Value = 361.0625
$ws.Range("M33").Value = -132.0625
$ws.Range("H40").Value = 6645.4546
$ws.Range("I40").Value = 11900.25
$ws.Range("K40").Value = 11900.25
$ws.Range("M40").Value = -11725.25
$ws.Range("H70").Value = 74395.37
$ws.Range("I70").Value = 114883.21
$ws.Range("J70").Value = 3541.625
$ws.Range("K70").Value = 344649.63
$ws.Range("L70").Value = 10624.875
$ws.Range("M70").Value = -344379.63
$ws.Range("N70").Value = -11164.875
$ws.Range("H73").Value = 74395.37
$ws.Range("I73").Value = 114883.21
$ws.Range("J73").Value = 3541.625
$ws.Range("K73").Value = 344649.63
$ws.Range("L73").Value = 10624.875
$ws.Range("M73").Value = -343713.63
$ws.Range("N73").Value = -12496.875
$ws.Range("H87").Value = 27333.334
$ws.Range("J87").Value = 27333.334
$ws.Range("L87").Value = 27333.334
$ws.Range("N87").Value = -29829.334
$ws.Range("H90").Value = 27333.334
$ws.Range("J90").Value = 27333.334
$ws.Range("L90").Value = 82000.00199999999
$ws.Range("N90").Value = -94480.00199999999
$ws.Range("H99").Value = 609.5714
$ws.Range("I99").Value = 452.6
$ws.Range("J99").Value = 1002
$ws.Range("K99").Value = 1357.8
$ws.Range("L99").Value = 3006
$ws.Range("M99").Value = 140.1999999999998
$ws.Range("N99").Value = -6002
$ws.Range("H106").Value = 3907
$ws.Range("I106").Value = 3907
$ws.Range("K106").Value = 3907
$ws.Range("M106").Value = -3276
$ws.Range("H107").Value = 8333
$ws.Range("I107").Value = 10000
$ws.Range("J107").Value = 6666
$ws.Range("K107").Value = 10000
$ws.Range("L107").Value = 6666
$ws.Range("M107").Value = -8080
$ws.Range("N107").Value = -10506
$ws.Range("H132").Value = 2663.9333
$ws.Range("J132").Value = 4009.5
$ws.Range("L132").Value = 12028.5
$ws.Range("N132").Value = -17088.5
$ws.Range("H137").Value = 3586
$ws.Range("I137").Value = 2646.8
$ws.Range("J137").Value = 4927.7144
$ws.Range("K137").Value = 7940.400000000001
$ws.Range("L137").Value = 14783.1432
$ws.Range("M137").Value = -5390.400000000001
$ws.Range("N137").Value = -19883.1432

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4122.1113
$ws.Range("I32").Value = 2201
$ws.Range("J32").Value = 16609.334
$ws.Range("K32").Value = 2201
$ws.Range("L32").Value = 16609.334
$ws.Range("M32").Value = -1914
$ws.Range("N32").Value = -17183.334
$ws.Range("H61").Value = 4817.1816
$ws.Range("I61").Value = 4388.5
$ws.Range("K61").Value = 4388.5
$ws.Range("M61").Value = -4176.5
$ws.Range("H74").Value = 1576.4546
$ws.Range("I74").Value = 1332.4286
$ws.Range("J74").Value = 2003.5
$ws.Range("K74").Value = 1332.4286
$ws.Range("L74").Value = 2003.5
$ws.Range("M74").Value = -458.4286
$ws.Range("N74").Value = -3751.5
$ws.Range("H77").Value = 1576.4546
$ws.Range("I77").Value = 1332.4286
$ws.Range("J77").Value = 2003.5
$ws.Range("K77").Value = 6662.143
$ws.Range("L77").Value = 10017.5
$ws.Range("M77").Value = -2294.143
$ws.Range("N77").Value = -18753.5
$ws.Range("H80").Value = 19333.334
$ws.Range("J80").Value = 40000
$ws.Range("L80").Value = 40000
$ws.Range("N80").Value = -41996
$ws.Range("H83").Value = 19333.334
$ws.Range("J83").Value = 40000
$ws.Range("L83").Value = 120000
$ws.Range("N83").Value = -129984
$ws.Range("H104").Value = 73266.664
$ws.Range("J104").Value = 73266.664
$ws.Range("L104").Value = 73266.664
$ws.Range("N104").Value = -80254.664
$ws.Range("H132").Value = 2316.25
$ws.Range("J132").Value = 2785.7273
$ws.Range("L132").Value = 8357.1819
$ws.Range("N132").Value = -13417.1819
$ws.Range("H136").Value = 4817.1816
$ws.Range("I136").Value = 4388.5
$ws.Range("K136").Value = 13165.5
$ws.Range("M136").Value = -10615.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 1801.1666
$ws.Range("J80").Value = 2129.8462
$ws.Range("L80").Value = 2129.8462
$ws.Range("N80").Value = -4125.8462
$ws.Range("H82").Value = 17322
$ws.Range("J82").Value = 40000
$ws.Range("L82").Value = 40000
$ws.Range("N82").Value = -40766
$ws.Range("H83").Value = 1801.1666
$ws.Range("J83").Value = 2129.8462
$ws.Range("L83").Value = 10649.231
$ws.Range("N83").Value = -20633.231
$ws.Range("H85").Value = 17322
$ws.Range("J85").Value = 40000
$ws.Range("L85").Value = 40000
$ws.Range("N85").Value = -42652
$ws.Range("H94").Value = 784.0769
$ws.Range("I94").Value = 808.5
$ws.Range("K94").Value = 808.5
$ws.Range("M94").Value = -357.5
$ws.Range("H107").Value = 2410.2222
$ws.Range("I107").Value = 2199.75
$ws.Range("J107").Value = 2578.6
$ws.Range("K107").Value = 2199.75
$ws.Range("L107").Value = 2578.6
$ws.Range("M107").Value = -279.75
$ws.Range("N107").Value = -6418.6
$ws.Range("H134").Value = 2615.926
$ws.Range("I134").Value = 2193.0833
$ws.Range("J134").Value = 5998.6665
$ws.Range("K134").Value = 6579.249899999999
$ws.Range("L134").Value = 17995.9995
$ws.Range("M134").Value = -4044.249899999999
$ws.Range("N134").Value = -23065.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3481.25
$ws.Range("I31").Value = 3176.2144
$ws.Range("K31").Value = 3176.2144
$ws.Range("M31").Value = -2881.2144
$ws.Range("H34").Value = 3481.25
$ws.Range("I34").Value = 3176.2144
$ws.Range("K34").Value = 3176.2144
$ws.Range("M34").Value = -2974.2144
$ws.Range("H68").Value = 88899
$ws.Range("J68").Value = 88899
$ws.Range("L68").Value = 88899
$ws.Range("N68").Value = -90397
$ws.Range("H71").Value = 88899
$ws.Range("J71").Value = 88899
$ws.Range("L71").Value = 266697
$ws.Range("N71").Value = -274185
$ws.Range("H74").Value = 39977.92
$ws.Range("J74").Value = 39977.92
$ws.Range("L74").Value = 39977.92
$ws.Range("N74").Value = -41725.92
$ws.Range("H77").Value = 39977.92
$ws.Range("J77").Value = 39977.92
$ws.Range("L77").Value = 119933.76
$ws.Range("N77").Value = -128669.76
$ws.Range("H99").Value = 12666.5
$ws.Range("I99").Value = 6162.3335
$ws.Range("J99").Value = 20471.5
$ws.Range("K99").Value = 6162.3335
$ws.Range("L99").Value = 20471.5
$ws.Range("M99").Value = -4664.3335
$ws.Range("N99").Value = -23467.5
$ws.Range("H126").Value = 12666.5
$ws.Range("I126").Value = 6162.3335
$ws.Range("J126").Value = 20471.5
$ws.Range("K126").Value = 18487.0005
$ws.Range("L126").Value = 61414.5
$ws.Range("M126").Value = -16017.0005
$ws.Range("N126").Value = -66354.5
$ws.Range("H134").Value = 2477.0454
$ws.Range("I134").Value = 2422.1333
$ws.Range("K134").Value = 7266.3999
$ws.Range("M134").Value = -4731.3999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 716
$ws.Range("J122").Value = 773.3333
$ws.Range("L122").Value = 6959.9997
$ws.Range("N122").Value = -11859.9997
$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("M140").Value = $null

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2586.739
$ws.Range("I102").Value = 2295.375
$ws.Range("K102").Value = 2295.375
$ws.Range("M102").Value = -673.375
$ws.Range("H126").Value = 4718.067
$ws.Range("I126").Value = 4001
$ws.Range("J126").Value = 4897.3335
$ws.Range("K126").Value = 12003
$ws.Range("L126").Value = 14692.0005
$ws.Range("M126").Value = -9533
$ws.Range("N126").Value = -19632.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5033.048
$ws.Range("J7").Value = 5579.2
$ws.Range("L7").Value = 5579.2
$ws.Range("N7").Value = -5803.2
$ws.Range("H40").Value = 5603.625
$ws.Range("J40").Value = 5741.727
$ws.Range("L40").Value = 5741.727
$ws.Range("N40").Value = -6013.727
$ws.Range("H68").Value = 2632.0667
$ws.Range("I68").Value = 2125.4546
$ws.Range("J68").Value = 4025.25
$ws.Range("K68").Value = 2125.4546
$ws.Range("L68").Value = 4025.25
$ws.Range("M68").Value = -1376.4546
$ws.Range("N68").Value = -5523.25
$ws.Range("H71").Value = 2632.0667
$ws.Range("I71").Value = 2125.4546
$ws.Range("J71").Value = 4025.25
$ws.Range("K71").Value = 10627.273
$ws.Range("L71").Value = 20126.25
$ws.Range("M71").Value = -6883.273000000001
$ws.Range("N71").Value = -27614.25
$ws.Range("H122").Value = 4899.7144
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 4899.7144
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 14699.1432
$ws.Range("M122").Value = $null
$ws.Range("N122").Value = -19599.1432
$ws.Range("H126").Value = 5033.048
$ws.Range("J126").Value = 5579.2
$ws.Range("L126").Value = 16737.6
$ws.Range("N126").Value = -21677.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7200.8
$ws.Range("I62").Value = 8000
$ws.Range("J62").Value = 7001
$ws.Range("K62").Value = 8000
$ws.Range("L62").Value = 7001
$ws.Range("M62").Value = -7376
$ws.Range("N62").Value = -8249
$ws.Range("H65").Value = 7200.8
$ws.Range("I65").Value = 8000
$ws.Range("J65").Value = 7001
$ws.Range("K65").Value = 40000
$ws.Range("L65").Value = 35005
$ws.Range("M65").Value = -36880
$ws.Range("N65").Value = -41245
$ws.Range("H100").Value = 323.83334
$ws.Range("I100").Value = 283.6
$ws.Range("J100").Value = 525
$ws.Range("K100").Value = 567.2
$ws.Range("L100").Value = 1050
$ws.Range("M100").Value = -26.20000000000005
$ws.Range("N100").Value = -2132
